$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 169 (pushes the existing rows 169:310 down to 170:311,
# growing the used range from A1:R310 to A1:R311).
$ws.Rows(169).Insert()

# Populate the newly inserted row with the new "Apio" price-report record.
$ws.Range("A169").Value = 5
$ws.Range("B169").Value = "Macroferia Regional de Talca"
$ws.Range("C169").Value = "Maule"
$ws.Range("D169").Value = 45090
$ws.Range("E169").Value = 7
$ws.Range("F169").Value = 100112017
$ws.Range("G169").Value = "Apio"
$ws.Range("H169").Value = "Americana (o)"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 70
$ws.Range("K169").Value = 5000
$ws.Range("L169").Value = 5000
$ws.Range("M169").Value = 5000
$ws.Range("N169").Value = "`$/docena de matas"
$ws.Range("O169").Value = "Provincia del Elquí"
$ws.Range("P169").Value = 833
$ws.Range("Q169").Value = 6
$ws.Range("R169").Value = "Hortaliza"
